$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.736399999999998
$ws.Range("D7").Value = -7.146399999999992
$ws.Range("B8").Value = 4.987100000000003
$ws.Range("A12").Value = -22.66730000000002
$ws.Range("B12").Value = 5.174300000000001
$ws.Range("B14").Value = 9.023600000000002
$ws.Range("D19").Value = -8.086699999999993
$ws.Range("E19").Value = 13.65460000000001
$ws.Range("D21").Value = -7.633499999999995
$ws.Range("B22").Value = 4.889000000000006
$ws.Range("D24").Value = -7.782699999999993
